$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 and Row 3 swap places (with a few value edits), and every row's
# "Forandrad" (C column) date changes from 45179 to 45180.
# ---------------------------------------------------------------------------

# New Row 2 content (was previously Row 3's case, "A 66380-2021")
$ws.Range("A2").Value = "A 66380-2021"
$ws.Range("B2").Value = 44518
$ws.Range("C2").Value = 45180
$ws.Range("D2").Value = "STOCKHOLMS LÄN"
$ws.Range("E2").Value = "HANINGE"
$ws.Range("G2").Value = 10.3
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 7
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 17
$ws.Range("R2").Value = "Ryl`r`nKnärot`r`nMotaggsvamp`r`nSpillkråka`r`nTallticka`r`nUllticka`r`nVedskivlav`r`nBlåmossa`r`nBronshjon`r`nFällmossa`r`nGrön sköldmossa`r`nGuldlockmossa`r`nJättesvampmal`r`nSårläka`r`nThomsons trägnagare`r`nVågbandad barkbock`r`nLopplummer"
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/artfynd/A 66380-2021.xlsx")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/kartor/A 66380-2021.png")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/knärot/A 66380-2021.png")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomål/A 66380-2021.docx")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomålsmail/A 66380-2021.docx")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsyn/A 66380-2021.docx")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsynsmail/A 66380-2021.docx")'

# New Row 3 content (was previously Row 2's case, "A 34417-2023") with a
# handful of updated counts and a trimmed species list.
$ws.Range("A3").Value = "A 34417-2023"
$ws.Range("B3").Value = 45139
$ws.Range("C3").Value = 45180
$ws.Range("D3").Value = "STOCKHOLMS LÄN"
$ws.Range("E3").Value = "HANINGE"
$ws.Range("G3").Value = 3.8
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 17
$ws.Range("R3").Value = "Barrviolspindling`r`nLeptoporus erubescens`r`nSpillkråka`r`nVedtrappmossa`r`nBronshjon`r`nDropptaggsvamp`r`nFällmossa`r`nGrön sköldmossa`r`nGuldlockmossa`r`nKornknutmossa`r`nRödgul trumpetsvamp`r`nStubbspretmossa`r`nSvavelriska`r`nSårläka`r`nVågbandad barkbock`r`nBlåsippa`r`nRevlummer"
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/artfynd/A 34417-2023.xlsx")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/kartor/A 34417-2023.png")'
$ws.Range("U3").ClearContents()
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomål/A 34417-2023.docx")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/klagomålsmail/A 34417-2023.docx")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsyn/A 34417-2023.docx")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_HANINGE/tillsynsmail/A 34417-2023.docx")'

# Keep the original fixed row heights (15pt) instead of Excel's automatic
# "best fit" height that kicks in after the multi-line text is rewritten.
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15

# ---------------------------------------------------------------------------
# Every remaining data row (4 through 91) only has its "Forandrad" date
# (column C) bumped from 45179 to 45180.
# ---------------------------------------------------------------------------
for ($r = 4; $r -le 91; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value = 45180
    }
}
